$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Swap the full data (columns B..AB) between each of these row pairs.
#    Column A (the running "id" index) and the row position itself stay put;
#    everything else (match id, teams, odds, results, ...) is exchanged
#    between the two rows of the pair.
# ---------------------------------------------------------------------------
$pairs = @(
    @(4,5),
    @(22,23),
    @(54,55),
    @(58,59),
    @(89,90),
    @(91,92),
    @(108,109),
    @(135,136),
    @(151,152),
    @(161,162),
    @(183,184)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    for ($c = 2; $c -le 28; $c++) {
        $v1 = $ws.Cells.Item($r1, $c).Value2
        $v2 = $ws.Cells.Item($r2, $c).Value2
        $ws.Cells.Item($r1, $c).Value2 = $v2
        $ws.Cells.Item($r2, $c).Value2 = $v1
    }
}

# ---------------------------------------------------------------------------
# 2) Append the three new match rows at the bottom of the sheet (196-198).
#    Copy the number formats for columns A (bold/bordered id) and D (date)
#    from the last pre-existing row so the new rows look like the rest.
# ---------------------------------------------------------------------------
$newRows = @(196, 197, 198)
foreach ($nr in $newRows) {
    $ws.Range("A195").Copy()
    $ws.Range("A$nr").PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    $ws.Range("D195").Copy()
    $ws.Range("D$nr").PasteSpecial(-4122)
    $excel.CutCopyMode = $false
}

function Set-Row($r, $vals) {
    $ws.Cells.Item($r, 1).Value2  = $vals.A
    $ws.Cells.Item($r, 2).Value2  = $vals.B
    $ws.Cells.Item($r, 3).Value2  = $vals.C
    $ws.Cells.Item($r, 4).Value2  = $vals.D
    $ws.Cells.Item($r, 5).Value2  = $vals.E
    $ws.Cells.Item($r, 6).Value2  = $vals.F
    $ws.Cells.Item($r, 7).Value2  = $vals.G
    $ws.Cells.Item($r, 8).Value2  = $vals.H
    $ws.Cells.Item($r, 9).Value2  = $vals.I
    $ws.Cells.Item($r, 10).Value2 = $vals.J
    $ws.Cells.Item($r, 11).Value2 = $vals.K
    $ws.Cells.Item($r, 12).Value2 = $vals.L
    $ws.Cells.Item($r, 13).Value2 = $vals.M
    $ws.Cells.Item($r, 14).Value2 = $vals.N
    $ws.Cells.Item($r, 15).Value2 = $vals.O
    $ws.Cells.Item($r, 16).Value2 = $vals.P
    $ws.Cells.Item($r, 17).Value2 = $vals.Q
    $ws.Cells.Item($r, 18).Value2 = $vals.R
    $ws.Cells.Item($r, 19).Value2 = $vals.S
    $ws.Cells.Item($r, 20).Value2 = $vals.T
    $ws.Cells.Item($r, 21).Value2 = $vals.U
    $ws.Cells.Item($r, 22).Value2 = $vals.V
    $ws.Cells.Item($r, 23).Value2 = $vals.W
    $ws.Cells.Item($r, 24).Value2 = $vals.X
    $ws.Cells.Item($r, 25).Value2 = $vals.Y
    $ws.Cells.Item($r, 26).Value2 = $vals.Z
    $ws.Cells.Item($r, 27).Value2 = $vals.AA
    $ws.Cells.Item($r, 28).Value2 = $vals.AB
}

$row196 = @{
    A = 194; B = 8192761; C = "Iraq League"; D = 45421.40625
    E = "Al Hudod"; F = "Al Quwa Al Jawiya"
    G = 0; H = 3; I = "A"
    J = 6; K = 3.6; L = 1.5
    M = 7; N = 3.75; O = 1.4
    P = 1.25; Q = 1.8; R = 2
    S = 2; T = 1.775; U = 2.025
    V = -1; W = -1; X = 0.3999999999999999
    Y = -1; Z = 1; AA = 0.7749999999999999
    AB = -1
}

$row197 = @{
    A = 195; B = 8192762; C = "Iraq League"; D = 45421.5
    E = "Karbalaa FC"; F = "Duhok"
    G = 1; H = 1; I = "D"
    J = 2.5; K = 2.75; L = 2.875
    M = 2.5; N = 2.875; O = 2.75
    P = 0; Q = 1.8; R = 2
    S = 1.75; T = 1.8; U = 2
    V = -1; W = 1.875; X = -1
    Y = 0; Z = 0; AA = 0.4
    AB = -0.5
}

$row198 = @{
    A = 196; B = 8197700; C = "Iraq League"; D = 45421.5
    E = "Amanat Baghdad"; F = "Zakho"
    G = 0; H = 0; I = "D"
    J = 4.75; K = 3.25; L = 1.666
    M = 3.5; N = 3.1; O = 1.95
    P = 0.5; Q = 1.775; R = 2.025
    S = 2; T = 1.975; U = 1.825
    V = -1; W = 2.1; X = -1
    Y = 0.7749999999999999; Z = -1; AA = -1
    AB = 0.825
}

Set-Row 196 $row196
Set-Row 197 $row197
Set-Row 198 $row198
